$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'245.38"
$ws.Cells.Item(2,5).Value = "'0.83%"
$ws.Cells.Item(2,7).Value = "'5"

$ws.Cells.Item(3,4).Value = "'29.33"
$ws.Cells.Item(3,5).Value = "'-2.20%"
$ws.Cells.Item(3,7).Value = "'5"

$ws.Cells.Item(4,4).Value = "'5.149"
$ws.Cells.Item(4,5).Value = "'0.28%"
$ws.Cells.Item(4,7).Value = "'5"

$ws.Cells.Item(5,4).Value = "'0.05775"
$ws.Cells.Item(5,5).Value = "'1.89%"
$ws.Cells.Item(5,7).Value = "'5"

$ws.Cells.Item(6,5).Value = "'1.40%"
$ws.Cells.Item(6,7).Value = "'5"

$ws.Cells.Item(7,4).Value = "'3.177"
$ws.Cells.Item(7,5).Value = "'5.15%"
$ws.Cells.Item(7,7).Value = "'5"

$ws.Cells.Item(8,4).Value = "'0.8573"
$ws.Cells.Item(8,5).Value = "'2.10%"
$ws.Cells.Item(8,7).Value = "'5"

$ws.Cells.Item(9,4).Value = "'0.8574"
$ws.Cells.Item(9,5).Value = "'-0.47%"
$ws.Cells.Item(9,7).Value = "'5"

$ws.Cells.Item(10,5).Value = "'2.12%"
$ws.Cells.Item(10,7).Value = "'5"

$ws.Cells.Item(11,4).Value = "'0.07031"
$ws.Cells.Item(11,5).Value = "'1.78%"
$ws.Cells.Item(11,7).Value = "'5"

$ws.Cells.Item(12,4).Value = "'0.03088"
$ws.Cells.Item(12,5).Value = "'7.82%"
$ws.Cells.Item(12,7).Value = "'5"

$ws.Cells.Item(13,4).Value = "'0.09372"
$ws.Cells.Item(13,5).Value = "'0.01%"
$ws.Cells.Item(13,7).Value = "'5"

$ws.Cells.Item(14,4).Value = "'0.001525"
$ws.Cells.Item(14,5).Value = "'-0.17%"
$ws.Cells.Item(14,7).Value = "'5"

$ws.Cells.Item(15,4).Value = "'0.0006020"
$ws.Cells.Item(15,5).Value = "'0.70%"
$ws.Cells.Item(15,7).Value = "'5"

$ws.Cells.Item(16,4).Value = "'0.006012"
$ws.Cells.Item(16,5).Value = "'0.48%"
$ws.Cells.Item(16,7).Value = "'5"

$ws.Cells.Item(17,5).Value = "'-0.86%"
$ws.Cells.Item(17,7).Value = "'5"

$ws.Cells.Item(18,4).Value = "'2.166"
$ws.Cells.Item(18,5).Value = "'1.80%"
$ws.Cells.Item(18,7).Value = "'5"

$ws.Cells.Item(19,4).Value = "'0.3203"
$ws.Cells.Item(19,5).Value = "'1.64%"
$ws.Cells.Item(19,7).Value = "'5"

$ws.Cells.Item(20,4).Value = "'0.03300"
$ws.Cells.Item(20,5).Value = "'1.05%"
$ws.Cells.Item(20,7).Value = "'5"

$ws.Cells.Item(21,4).Value = "'0.1282"
$ws.Cells.Item(21,5).Value = "'-1.08%"
$ws.Cells.Item(21,7).Value = "'5"

$ws.Cells.Item(22,4).Value = "'3.178"
$ws.Cells.Item(22,5).Value = "'-12.45%"
$ws.Cells.Item(22,7).Value = "'5"

$ws.Cells.Item(23,4).Value = "'0.04145"
$ws.Cells.Item(23,5).Value = "'-0.52%"
$ws.Cells.Item(23,7).Value = "'5"

$ws.Cells.Item(24,5).Value = "'1.92%"
$ws.Cells.Item(24,7).Value = "'5"

$ws.Cells.Item(25,5).Value = "'1.42%"
$ws.Cells.Item(25,7).Value = "'5"

$ws.Cells.Item(26,4).Value = "'0.004130"
$ws.Cells.Item(26,7).Value = "'5"

$ws.Cells.Item(27,5).Value = "'2.58%"
$ws.Cells.Item(27,7).Value = "'5"

$ws.Cells.Item(28,5).Value = "'3.42%"
$ws.Cells.Item(28,7).Value = "'5"

$ws.Cells.Item(29,7).Value = "'5"

$ws.Cells.Item(30,7).Value = "'5"

$ws.Cells.Item(31,7).Value = "'5"

$ws.Cells.Item(32,7).Value = "'5"

$ws.Cells.Item(33,7).Value = "'5"

$ws.Cells.Item(34,7).Value = "'5"

$ws.Cells.Item(35,7).Value = "'5"

$ws.Cells.Item(36,7).Value = "'5"

$ws.Cells.Item(37,7).Value = "'5"

$ws.Cells.Item(38,7).Value = "'5"

$ws.Cells.Item(39,7).Value = "'5"

$ws.Cells.Item(40,4).Value = "'0.03725"
$ws.Cells.Item(40,5).Value = "'0.26%"
$ws.Cells.Item(40,7).Value = "'5"

$ws.Cells.Item(41,2).Value = "BKEXToken"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(41,4).Value = "'0.1070"
$ws.Cells.Item(41,5).Value = "'1.09%"
$ws.Cells.Item(41,7).Value = "'5"

$ws.Cells.Item(42,2).Value = "CEJI"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(42,4).Value = "'0.002450"
$ws.Cells.Item(42,5).Value = "'6.04%"
$ws.Cells.Item(42,7).Value = "'5"

$ws.Cells.Item(43,2).Value = "KickToken"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(43,4).Value = "'0.003522"
$ws.Cells.Item(43,5).Value = "'-33.90%"
$ws.Cells.Item(43,7).Value = "'5"

$ws.Cells.Item(44,4).Value = "'0.009172"
$ws.Cells.Item(44,5).Value = "'-5.30%"
$ws.Cells.Item(44,7).Value = "'5"

$ws.Cells.Item(45,4).Value = "'0.00005283"
$ws.Cells.Item(45,5).Value = "'3.54%"
$ws.Cells.Item(45,7).Value = "'5"

$ws.Cells.Item(46,7).Value = "'5"

$ws.Cells.Item(47,4).Value = "'0.05800"
$ws.Cells.Item(47,5).Value = "'-41.98%"
$ws.Cells.Item(47,7).Value = "'5"

$ws.Cells.Item(48,5).Value = "'-20.19%"
$ws.Cells.Item(48,7).Value = "'5"

$ws.Cells.Item(49,7).Value = "'5"

$ws.Cells.Item(50,7).Value = "'5"

$ws.Cells.Item(51,7).Value = "'5"
